$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-numeric-looking price cells to remain Text, matching source formatting
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = '59.310.92'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '2.514.95'
$ws.Range("E3").Value = '  -0.21%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '536.52'
$ws.Range("E5").Value = '  -0.87%  '
$ws.Range("D6").Value = '139.40'
$ws.Range("E6").Value = '  -3.71%  '
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E8").Value = '  -1.48%  '
$ws.Range("D9").Value = '2.517.37'
$ws.Range("E9").Value = '  -0.97%  '
$ws.Range("E10").Value = '  +0.45%  '
$ws.Range("E11").Value = '  +1.48%  '
$ws.Range("D12").Value = '5.42'
$ws.Range("E12").Value = '  -3.32%  '
$ws.Range("E13").Value = '  -0.47%  '
$ws.Range("D14").Value = '2.961.28'
$ws.Range("E14").Value = '  -0.13%  '
$ws.Range("D15").Value = '23.47'
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").Value = '59.209.09'
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").Value = '2.519.68'
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").Value = '11.16'
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("E20").Value = '  +0.70%  '
$ws.Range("D21").Value = '324.91'
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '5.81'
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("D24").Value = '63.64'
$ws.Range("E24").Value = '  +2.60%  '
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("D26").Value = '0.168'
$ws.Range("E26").Value = '  +1.67%  '
$ws.Range("D28").Value = '7.82'
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("D29").Value = '6.96'
$ws.Range("E29").Value = '  +4.05%  '
$ws.Range("D30").Value = '0.0₃0776'
$ws.Range("E30").Value = '  -0.59%  '
$ws.Range("E31").Value = '  -1.95%  '
$ws.Range("D32").Value = '166.29'
$ws.Range("E32").Value = '  +4.89%  '
$ws.Range("E33").Value = '  +0.14%  '
$ws.Range("D34").Value = '1.45'
$ws.Range("E34").Value = '  -2.78%  '
$ws.Range("E35").Value = '  -7.29%  '
$ws.Range("D36").Value = '18.52'
$ws.Range("E36").Value = '  -1.05%  '
$ws.Range("D37").Value = '4.27'
$ws.Range("E37").Value = '  -2.40%  '
$ws.Range("E38").Value = '  -1.48%  '
$ws.Range("D39").Value = '36.89'
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("E40").Value = '  -0.36%  '
$ws.Range("D41").Value = '0.814'
$ws.Range("E41").Value = '  -1.29%  '
$ws.Range("D42").Value = '5.25'
$ws.Range("E42").Value = '  -6.53%  '
$ws.Range("D43").Value = '279.71'
$ws.Range("E43").Value = '  -4.96%  '
$ws.Range("E45").Value = '  +0.70%  '
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("D48").Value = '122.87'
$ws.Range("E48").Value = '  +0.53%  '
$ws.Range("D49").Value = '0.0514'
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("E51").Value = '  -1.84%  '
